$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A new hourly observation for 2026/02/21 (Sat) was recorded between the
# existing 2026/02/21 03:00 row (851) and the first 2026/12/29 row (old
# row 851). Insert a fresh row at 851, pushing rows 851:892 down to
# 852:893 and growing the used range from D892 to D893.
$ws.Rows.Item(851).Insert()

# Column A holds date-like text (e.g. "2026/02/21") stored as plain
# strings in this sheet, not real dates. Assigning that text straight to
# .Value/.Formula makes Excel "smart"-convert it into a date serial
# (and stamp a date NumberFormat on the cell), which would not match the
# source data. Routing it through a literal-text formula and then
# collapsing the formula to its value via Copy + PasteSpecial (values
# only) keeps it a plain string with no extra formatting applied.
$ws.Cells.Item(851, 1).Formula = '="2026/02/21"'
$ws.Cells.Item(851, 1).Copy()
$ws.Cells.Item(851, 1).PasteSpecial(-4163) # xlPasteValues

$ws.Cells.Item(851, 2).Value = "土"
$ws.Cells.Item(851, 3).Value = 7
$ws.Cells.Item(851, 4).Value = 201

$excel.CutCopyMode = $false

$wb.Save()
